$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 132, shifting the existing rows 132-134 down to 133-135.
$ws.Rows.Item(132).Insert()

# Populate the newly inserted row 132 with the new weekly entry.
$ws.Cells.Item(132, 1).Value = 4
$ws.Cells.Item(132, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(132, 3).Value = "Los Lagos"
$ws.Cells.Item(132, 4).Value = 44448
$ws.Cells.Item(132, 5).Value = 10
$ws.Cells.Item(132, 6).Value = "Fruta"
$ws.Cells.Item(132, 7).Value = 100108
$ws.Cells.Item(132, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(132, 9).Value = 100108005
$ws.Cells.Item(132, 10).Value = "Piña"
$ws.Cells.Item(132, 11).Value = "Caramelo"
$ws.Cells.Item(132, 12).Value = "Primera"
$ws.Cells.Item(132, 13).Value = 80
$ws.Cells.Item(132, 14).Value = 24000
$ws.Cells.Item(132, 15).Value = 24000
$ws.Cells.Item(132, 16).Value = 24000
$ws.Cells.Item(132, 17).Value = "$/caja 12 unidades"
$ws.Cells.Item(132, 18).Value = "Ecuador"
$ws.Cells.Item(132, 19).Value = 2000
$ws.Cells.Item(132, 20).Value = 12
